# 55927: "handle date formulas too" - add a second XML-mapped date cell
# whose value comes from a DATE() formula rather than a literal value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed A2 from A1 (copy/paste) so it inherits A1's cell style (the date
# number format), then overwrite it with the date formula. A straight
# Formula assignment on a blank cell would pick up the default (General)
# style instead of the date format used by A1.
$ws.Range("A1").Copy($ws.Range("A2")) | Out-Null
$ws.Range("A2").Formula = "=DATE(2012,2,16)"

# Move the active selection on to A3, matching the state left behind after
# typing the formula into A2 and pressing Enter.
$ws.Range("A3").Select() | Out-Null
